$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "test_controle_adc_10bits" (3rd sheet): append a new results table
#    ("tension appliquée" / "tension mesurée") below the existing one.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B18").Value = "canal"
$ws3.Range("C18").Value = "tension appliquée"
$ws3.Range("D18").Value = "tension mesurée"

$rows3 = @(
    @("1 V",      "0.9989 V"),
    @("-1V",      "-1.0029 V"),
    @("1.5 V",    "1.493 V"),
    @("-1.5 V",   "-1.511 V"),
    @("0 V",      "-0.008 ou 0 ou 0.008 V"),
    @("0.5 V",    "0.5045 V"),
    @("-0.5 V",   "-0.4965 V"),
    @("0.7546 V", "0.7587 V"),
    @("-1.045 V", "-1.037 V"),
    @("1.845 V",  "1,851 V")
)

$r = 19
foreach ($row in $rows3) {
    $ws3.Range("B$r").Value = 0
    $ws3.Range("C$r").Value = $row[0]
    $ws3.Range("D$r").Value = $row[1]
    $r++
}

# Match formatting (centered, no border) used by the rest of the sheet.
$ws3.Range("B18:F28").HorizontalAlignment = -4108

$ws3.Range("E28").Select() | Out-Null

# ---------------------------------------------------------------------------
# 1b) Leftover cursor positions on the other pre-existing sheets.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F16").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) New sheet "test_adc_dac_modele_rc" at the end of the workbook, containing
#    the RC model test results.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "test_adc_dac_modele_rc"

$ws5.Range("B3").Value = "consigne DAC"
$ws5.Range("C3").Value = "courant mesuré"
$ws5.Range("D3").Value = "tension test1"
$ws5.Range("E3").Value = "tension test2"

$rows5 = @(
    @("1,8958 V",  "-0.74 V"),
    @("1,2626 V",  "-0.43 V"),
    @("0,6317 V",  "-0.65 V "),
    @("0 V",       "-0.018 V "),
    @(-0.507,      "0.4925 V"),
    @("-1,044 V",  "0.991 V"),
    @("-1,5209 V", "1.233 V")
)

$r = 4
foreach ($row in $rows5) {
    $ws5.Range("B$r").Value = $row[0]
    $ws5.Range("D$r").Value = $row[1]
    $r++
}

# Empty but formatted "courant mesuré" column.
$ws5.Range("B3:D10").HorizontalAlignment = -4108

$ws5.Range("C4").Select() | Out-Null
